$d = $word.ActiveDocument

$replacements = @(
    @("499÷9=55, 4", "332÷6=55, 2"),
    @("122÷9=13, 5", "575÷5=115, 0"),
    @("354÷7=50, 4", "567÷4=141, 3"),
    @("427÷3=142, 1", "646÷8=80, 6"),
    @("433÷6=72, 1", "172÷3=57, 1"),
    @("514÷7=73, 3", "361÷6=60, 1"),
    @("444÷7=63, 3", "536÷4=134, 0"),
    @("849÷9=94, 3", "523÷5=104, 3"),
    @("835÷4=208, 3", "703÷4=175, 3"),
    @("937÷6=156, 1", "966÷3=322, 0"),
    @("356÷8=44, 4", "839÷6=139, 5"),
    @("245÷9=27, 2", "677÷8=84, 5"),
    @("143÷9=15, 8", "693÷2=346, 1"),
    @("216÷3=72, 0", "823÷9=91, 4"),
    @("609÷3=203, 0", "633÷3=211, 0"),
    @("255÷5=51, 0", "283÷2=141, 1"),
    @("143÷7=20, 3", "406÷4=101, 2"),
    @("799÷2=399, 1", "734÷6=122, 2"),
    @("184÷9=20, 4", "460÷3=153, 1"),
    @("864÷5=172, 4", "817÷4=204, 1"),
    @("329÷9=36, 5", "687÷2=343, 1"),
    @("846÷8=105, 6", "726÷9=80, 6"),
    @("157÷2=78, 1", "972÷4=243, 0"),
    @("823÷5=164, 3", "567÷8=70, 7"),
    @("299÷7=42, 5", "547÷9=60, 7")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
